# Insert a new weekly record at row 55 ("Vega Central Mapocho de Santiago -
# Bruselas (repollito)" sheet). Inserting the row pushes the existing rows
# 55-116 down to 56-117 (dimension grows from A1:R116 to A1:R117), matching
# the rest of the data set which keeps growing downward with each new week.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(55).Insert()

$ws.Range("A55").Value = 9
$ws.Range("B55").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C55").Value = "Metropolitana"
$ws.Range("D55").Value = 45167
$ws.Range("E55").Value = 13
$ws.Range("F55").Value = 100112035
$ws.Range("G55").Value = "Bruselas (repollito)"
$ws.Range("H55").Value = "Sin especificar"
$ws.Range("I55").Value = "Primera"
$ws.Range("J55").Value = 52
$ws.Range("K55").Value = 16000
$ws.Range("L55").Value = 16000
$ws.Range("M55").Value = 16000
$ws.Range("N55").Value = "$/malla 15 kilos"
$ws.Range("O55").Value = "Provincia de Quillota"
$ws.Range("P55").Value = 1067
$ws.Range("Q55").Value = 15
$ws.Range("R55").Value = "Hortaliza"
